$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new time-tracking entry for row 21 (2020-08-03 as Excel serial date)
$ws.Range("A21").Value = 44046
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = "Oauthin toiminnan selvittämistä ja usereiden tallentamisen ja toiminnan suunnittelua"

# Row grows to two text lines, same as the other wrapped-text rows (19/20)
$ws.Rows.Item(21).RowHeight = 30

# Update the active selection to reflect where the user ended up working
$ws.Range("G21").Select()
